$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking values keep their exact
# original text representation (trailing zeros, "xx.xxx.xx"-style ids, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.907.43'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '1.859.72'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Value = '1.017'
$ws.Range("E4").Value = '  -1.88%  '
$ws.Range("D5").Value = '320.93'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").Value = '1.016'
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("D7").Value = '0.4341'
$ws.Range("E7").Value = '  -1.52%  '
$ws.Range("D8").Value = '0.3795'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.07444'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '0.8879'
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("D11").Value = '21.76'
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = '1.871.31'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '6.808'
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").Value = '5.505'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = '0.07119'
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").Value = '88.44'
$ws.Range("E16").Value = '  +5.23%  '
$ws.Range("D17").Value = '1.022'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").Value = '0.000009053'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("D20").Value = '15.49'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '27.919.47'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = '5.282'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").Value = '11.23'
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").Value = '2.076.86'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").Value = '2.031'
$ws.Range("E25").Value = '  +4.86%  '
$ws.Range("D26").Value = '157.21'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("D27").Value = '18.73'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").Value = '5.440'
$ws.Range("E28").Value = '  +2.14%  '
$ws.Range("D29").Value = '2.015'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").Value = '120.40'
$ws.Range("E30").Value = '  +2.28%  '
$ws.Range("D31").Value = '0.09022'
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("D32").Value = '1.242'
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("D33").Value = '0.7762'
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").Value = '4.595'
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").Value = '2.997'
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").Value = '1.017'
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("D37").Value = '1.146'
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("D38").Value = '0.01980'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").Value = '0.05320'
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").Value = '2.885'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("D41").Value = '0.5216'
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").Value = '7.001'
$ws.Range("E42").Value = '  +2.10%  '
$ws.Range("D43").Value = '0.1683'
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").Value = '8.802'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("D45").Value = '110.89'
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("D46").Value = '10.77'
$ws.Range("E46").Value = '  +1.19%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.4764'
$ws.Range("E47").Value = '  +1.56%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.717'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").Value = '1.017'
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06489'
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").Value = '1.897'
$ws.Range("E51").Value = '  +1.14%  '
